$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> column -> new value, derived from the diff.
$updates = @{
    17  = @{ B = 5;  D = 1; E = 0;   F = 2;   G = 2.5 }
    39  = @{ B = 12; C = 2; E = 0;   F = 4;   G = 3 }
    41  = @{ B = 5;  D = 1; E = 0;   F = 2;   G = 2.5 }
    66  = @{ B = 6;  G = 1.5 }
    88  = @{ B = 11; D = 2; E = 0.5; F = 3.5; G = 3.142857142857143 }
    93  = @{ B = 21; G = 5.25 }
    110 = @{ B = 11; D = 2; E = 0.5; F = 3.5; G = 3.142857142857143 }
    111 = @{ B = 14; D = 2; F = 4;   G = 3.5 }
    133 = @{ B = 3;  G = 0.75 }
    134 = @{ B = 11; D = 2; F = 4;   G = 2.75 }
    135 = @{ B = 14; D = 2; F = 4;   G = 3.5 }
    136 = @{ B = 11; D = 2; E = 0.5; F = 3.5; G = 3.142857142857143 }
    158 = @{ B = 10; C = 1; D = 3;   F = 4;   G = 2.5 }
    159 = @{ B = 13; C = 1; D = 1; E = 1; F = 3; G = 4.333333333333333 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
